$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Resistance Level -> Resistance Confirmed
$ws.Range("C3").Value = 6714.7900390625
$ws.Range("E3").Value = "Resistance Confirmed"
$ws.Range("F3").Value = 49.41015625

# Row 4: Resistance Level -> Resistance Confirmed
$ws.Range("C4").Value = 6839.72998046875
$ws.Range("E4").Value = "Resistance Confirmed"
$ws.Range("F4").Value = 52.89013671875

# Row 5: Resistance Level -> Resistance Confirmed
$ws.Range("C5").Value = 6848.33984375
$ws.Range("E5").Value = "Resistance Confirmed"
$ws.Range("F5").Value = 31.27001953125

# Row 16: Support Level -> Closing Price
$ws.Range("C16").Value = 6840.31982421875
$ws.Range("E16").Value = "Closing Price"
$ws.Range("F16").Value = 13.8203125

# Row 20: Support Level -> Support Confirmed
$ws.Range("C20").Value = 6767.27001953125
$ws.Range("E20").Value = "Support Confirmed"
$ws.Range("F20").Value = 37.72021484375

# Row 22: Resistance Level -> Resistance Confirmed
$ws.Range("C22").Value = 6909.740234375
$ws.Range("E22").Value = "Resistance Confirmed"
$ws.Range("F22").Value = 8.76025390625

# Row 23: Support Level -> Closing Price
$ws.Range("C23").Value = 6896.4501953125
$ws.Range("E23").Value = "Closing Price"
$ws.Range("F23").Value = 3.69970703125

# Row 26: Support Level -> Closing Price
$ws.Range("C26").Value = 6944.97021484375
$ws.Range("E26").Value = "Closing Price"
$ws.Range("F26").Value = 27.1298828125
